$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "67.285.50"
$ws.Range("E2").Value = "  +4.53%  "

$ws.Range("D3").Value = "3.454.32"
$ws.Range("E3").Value = "  +2.89%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.44"
$ws.Range("E5").Value = "  +2.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.20"
$ws.Range("E6").Value = "  +6.07%  "

$ws.Range("E7").Value = "  +2.27%  "

$ws.Range("D8").Value = "3.447.18"
$ws.Range("E8").Value = "  +2.96%  "

$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("E10").Value = "  +8.26%  "

$ws.Range("E11").Value = "  +2.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.56"
$ws.Range("E12").Value = "  +3.74%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000280"
$ws.Range("E13").Value = "  +2.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.37"
$ws.Range("E14").Value = "  +3.82%  "

$ws.Range("D15").Value = "4.005.72"
$ws.Range("E15").Value = "  +2.77%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.462.32"
$ws.Range("E16").Value = "  +3.10%  "

$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.54"
$ws.Range("E17").Value = "  +2.16%  "

$ws.Range("D18").Value = "67.239.16"
$ws.Range("E18").Value = "  +4.34%  "

$ws.Range("E19").Value = "  +1.20%  "

$ws.Range("E20").Value = "  +2.81%  "

$ws.Range("E21").Value = "  +2.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "480.85"
$ws.Range("E22").Value = "  +5.20%  "

$ws.Range("E23").Value = "  +1.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.10"
$ws.Range("E24").Value = "  +11.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.18"
$ws.Range("E25").Value = "  +2.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.75"
$ws.Range("E26").Value = "  +4.62%  "

$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.91"
$ws.Range("E28").Value = "  +1.85%  "

$ws.Range("E29").Value = "  +2.88%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.53"
$ws.Range("E30").Value = "  +3.85%  "

$ws.Range("E31").Value = "  +4.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "601.89"
$ws.Range("E32").Value = "  +5.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.61"
$ws.Range("E33").Value = "  +1.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "62.95"
$ws.Range("E34").Value = "  +3.37%  "

$ws.Range("E35").Value = "  +2.47%  "

$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("E37").Value = "  +6.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.66"
$ws.Range("E38").Value = "  +1.05%  "

$ws.Range("E39").Value = "  +5.92%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.388"
$ws.Range("E40").Value = "  +5.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.52"
$ws.Range("E41").Value = "  +3.75%  "

$ws.Range("D42").Value = "3.152.30"
$ws.Range("E42").Value = "  +2.77%  "

$ws.Range("E43").Value = "  +4.28%  "

$ws.Range("E44").Value = "  +5.93%  "

$ws.Range("E45").Value = "  +3.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.81"
$ws.Range("E46").Value = "  +22.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.26"
$ws.Range("E47").Value = "  +3.62%  "

$ws.Range("E48").Value = "  +1.29%  "

$ws.Range("E49").Value = "  +7.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "142.19"
$ws.Range("E51").Value = "  +2.95%  "
